# Q3 Update - 2025
# Applies the data refresh to the "fromCSV" sheet:
#  - short-url column (B) value changes for every data row
#  - a handful of refugee/asylum-seeker counts are corrected (rows 509-513, 516)
#  - the country-of-origin details on rows 516-517 are updated
#  - the last data row (518) is removed (the table now ends at row 517)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues - used to push a value onto a cell while keeping that
# cell's existing number format/style (so numeric-looking strings like
# "12" stay stored as text, matching the rest of the column).
$xlPasteValues = -4163

function Set-TextValue($cellAddress, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellAddress).PasteSpecial($xlPasteValues)
    $scratch.Clear()
}

# 1. short-url column: "0fmXSZ" -> "74NUgW" for every data row (2-518)
for ($r = 2; $r -le 518; $r++) {
    $ws.Cells.Item($r, 2).Value = "74NUgW"
}

# 2. Refugee / asylum-seeker figure corrections
Set-TextValue "O509" "12"
Set-TextValue "N510" "20"
Set-TextValue "O511" "9"
Set-TextValue "N512" "55"
Set-TextValue "O512" "5"
Set-TextValue "N513" "61"

# 3. Row 516 now reports on Sudan instead of Sao Tome and Principe
Set-TextValue "F516" "177"
$ws.Range("G516").Value = "Sudan"
$ws.Range("H516").Value = "SUD"
$ws.Range("I516").Value = "SDN"
Set-TextValue "N516" "11"

# 4. Row 517 now reports on Togo instead of Sudan
Set-TextValue "F517" "192"
$ws.Range("G517").Value = "Togo"
$ws.Range("H517").Value = "TOG"
$ws.Range("I517").Value = "TGO"

# 5. Drop the final data row - the table now ends at row 517
$ws.Rows.Item(518).Delete()
